$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new "Build command" section: insert two blank rows
# above the old row 22 ("Merge to NCS." + the numbered readme steps),
# pushing that block down to rows 24-26.
$ws.Range("A22:A23").EntireRow.Insert()
$ws.Range("A22").Value = "Build command"

# Insert a "TODO" marker (bold, yellow fill) in column A of row 18,
# next to the "Replace in Wi-Fi and BLE coex sample." note.
$ws.Range("A18").Value = "TODO"
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").Interior.Color = 65535

# Default page setup (portrait, paper size 9 / A4), matching the
# printable-sheet settings saved with the rest of the edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the view scrolled/selected near the new content.
$ws.Range("E14").Select()
